$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Scottie Barnes", "PG,SG,SF,PF", "Toronto Raptors"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Nick Richards", "C", "Phoenix Suns"),
    @("Max Strus", "SG,SF", "Cleveland Cavaliers"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Isaiah Collier", "PG", "Utah Jazz"),
    @("Ja Morant", "PG", "Memphis Grizzlies"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Clint Capela", "C", "Atlanta Hawks")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
